$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.431.49"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.804.48"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.76"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.27%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.78%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "38.12"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +4.61%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -5.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0669"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.63%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.064.06"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -7.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.810.24"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.418.02"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.37"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.67"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.25"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0766"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.00"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -5.28%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.08"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.77%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.61"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.69"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.37"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.118"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.73"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.01%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.83"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.12%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0511"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.81"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.322.37"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.42%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.637"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -5.72%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.05"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.18%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.66%  "
$ws.Range("B39").Value = "HuobiToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.45"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.97%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.30"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -7.61%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.81"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.21"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.43%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "81.27"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.85%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.62"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0510"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.964.80"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -5.18%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "101.73"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0119"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -6.71%  "
